$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cells from the cryptos list refresh.
# NumberFormat is forced to Text ("@") before assignment so that
# numeric-looking strings (e.g. "207.24") are not auto-converted to
# floating point numbers, then restored to "General" afterwards so no
# spurious style/format diff is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.980.02'
$ws.Range("D2").NumberFormat = "General"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.560.94'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.56%  '
$ws.Range("E3").NumberFormat = "General"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("E4").NumberFormat = "General"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.24'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("E5").NumberFormat = "General"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.48%  '
$ws.Range("E6").NumberFormat = "General"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("E7").NumberFormat = "General"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.12'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.82%  '
$ws.Range("E8").NumberFormat = "General"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("E9").NumberFormat = "General"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0596'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.93%  '
$ws.Range("E10").NumberFormat = "General"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0858'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("E11").NumberFormat = "General"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.782.38'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("E12").NumberFormat = "General"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.545.49'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.25%  '
$ws.Range("E13").NumberFormat = "General"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.59%  '
$ws.Range("E14").NumberFormat = "General"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.86%  '
$ws.Range("E15").NumberFormat = "General"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.06'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.64%  '
$ws.Range("E16").NumberFormat = "General"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.978.35'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.28%  '
$ws.Range("E17").NumberFormat = "General"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '217.27'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("E18").NumberFormat = "General"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.06%  '
$ws.Range("E19").NumberFormat = "General"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.13%  '
$ws.Range("E20").NumberFormat = "General"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("E21").NumberFormat = "General"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.33%  '
$ws.Range("E22").NumberFormat = "General"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.55%  '
$ws.Range("E24").NumberFormat = "General"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.57'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.19%  '
$ws.Range("E25").NumberFormat = "General"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.37%  '
$ws.Range("E26").NumberFormat = "General"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.05'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.24%  '
$ws.Range("E27").NumberFormat = "General"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.37%  '
$ws.Range("E28").NumberFormat = "General"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("E29").NumberFormat = "General"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.67%  '
$ws.Range("E30").NumberFormat = "General"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.11'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.64%  '
$ws.Range("E31").NumberFormat = "General"

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.64%  '
$ws.Range("E32").NumberFormat = "General"

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("B33").NumberFormat = "General"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("C33").NumberFormat = "General"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.11'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.75%  '
$ws.Range("E33").NumberFormat = "General"

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'Maker'
$ws.Range("B34").NumberFormat = "General"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("C34").NumberFormat = "General"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.422.53'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("E34").NumberFormat = "General"

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.05%  '
$ws.Range("E35").NumberFormat = "General"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.05'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +9.67%  '
$ws.Range("E36").NumberFormat = "General"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.27%  '
$ws.Range("E37").NumberFormat = "General"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.73%  '
$ws.Range("E38").NumberFormat = "General"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.532'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.17%  '
$ws.Range("E39").NumberFormat = "General"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("E40").NumberFormat = "General"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("E41").NumberFormat = "General"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("E42").NumberFormat = "General"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.87%  '
$ws.Range("E43").NumberFormat = "General"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.90%  '
$ws.Range("E44").NumberFormat = "General"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.81'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.77%  '
$ws.Range("E45").NumberFormat = "General"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.29%  '
$ws.Range("E46").NumberFormat = "General"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.695.84'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.48%  '
$ws.Range("E47").NumberFormat = "General"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.33'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.34%  '
$ws.Range("E48").NumberFormat = "General"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.17%  '
$ws.Range("E49").NumberFormat = "General"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.36%  '
$ws.Range("E50").NumberFormat = "General"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.41%  '
$ws.Range("E51").NumberFormat = "General"
